# Update a few entries in the "Discussion" column (E) on the schedule.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "Scientific -> statistical question"
$ws.Range("E11").Value = "More on linear regression in R"
$ws.Range("E30").Value = "Scientific paper critique?"

# Leave the cursor where the author left it when they saved.
$ws.Range("C31").Select()
